# Add a "Meta description" paragraph right after the title (Heading1) and
# update the duplicated title/description block near the end of the document:
#   - drop the stray duplicate bold title paragraph
#   - turn the italic "Read our review..." paragraph into the DALLE image prompt

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the Heading1
#    title paragraph, matching the surrounding body paragraphs' run layout
#    (a leading empty run, a bold "Meta description" run, then the rest).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$insPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$metaXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Collapsed Castle Bonus Buy, a slot game with high RTP and volatility. Play for free and enjoy its medieval and fantasy theme.</w:t></w:r></w:p>
          <w:p></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insPoint.InsertXML($metaXml)

# InsertXML needs a trailing paragraph mark to avoid merging with the next
# paragraph, which leaves behind one extra empty paragraph - remove it.
# (paragraph 1 = title, paragraph 2 = new "Meta description" paragraph,
#  paragraph 3 = the leftover empty spacer paragraph)
$spacerPara = $d.Paragraphs.Item(3)
$spacerPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicate bold title paragraph
#    and rewrite the italic paragraph's text with the DALLE prompt.
# ---------------------------------------------------------------------------
$boldTitlePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$boldTitlePara.Range.Delete()

$descPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$descRange = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$descRange.Text = 'Prompt: DALLE, please create a cartoon-style feature image for the game "Collapsed Castle Bonus Buy" that features a happy Maya warrior with glasses. The image should incorporate elements of the medieval fantasy theme, such as a castle in the background or treasure in the foreground. Make sure the image is eye-catching and highlights the adventurous nature of the game.'
